$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.988.22"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "3.428.17"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").Value = "128.78"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +6.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.740"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.36%  "
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  +4.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "0.0000223"
$ws.Range("E12").Value = "  +48.55%  "
$ws.Range("D13").Value = "9.18"
$ws.Range("E13").Value = "  +9.39%  "
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("E15").Value = "  +7.90%  "
$ws.Range("D16").Value = "3.967.40"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "3.444.22"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  +6.84%  "
$ws.Range("E19").Value = "  +7.25%  "
$ws.Range("D20").Value = "61.987.39"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "450.39"
$ws.Range("E21").Value = "  +44.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "91.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.19%  "
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").Value = "13.02"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("D25").Value = "3.24"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").Value = "32.97"
$ws.Range("E26").Value = "  +10.91%  "
$ws.Range("D27").Value = "8.84"
$ws.Range("E27").Value = "  +8.52%  "
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "7.68"
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").Value = "12.01"
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("D33").Value = "43.14"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("E36").Value = "  +3.21%  "
$ws.Range("D37").Value = "54.44"
$ws.Range("E37").Value = "  +5.38%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("E40").Value = "  +7.62%  "
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("D43").Value = "142.28"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  +8.78%  "
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("E46").Value = "  +13.42%  "
$ws.Range("D47").Value = "16.69"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "22.32"
$ws.Range("E48").Value = "  +4.91%  "
$ws.Range("D49").Value = "2.13"
$ws.Range("E49").Value = "  +9.17%  "
$ws.Range("D50").Value = "3.775.68"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("E51").Value = "  +15.18%  "
